$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a header label in A1, reusing the same formatting as the other
# header cells (B1:W1) by copying their format only (keeps the existing
# style definition instead of creating a brand-new one).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1").Value = "Category"

# The category column (A2:A46) no longer carries the bold/bordered
# header style - reset it back to the default "Normal" style.
$ws.Range("A2:A46").Style = "Normal"
